$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testData")

# addDepartment test row: bump generated rowIndex value 3 -> 4
$ws.Range("D9").Value = "4"

# updateDepartment test row: rowIndex value 4 -> 2
$ws.Range("D24").Value = "2"

# Add the new "removeDepartment" automation scenario (rows 26-27),
# mirroring the existing "removeCompany" scenario (rows 20-21) for
# formatting, then relabel the text + restore the trailing blank
# marker cell in column F (matching the addDepartment/login-style rows).
$ws.Range("A20:E21").Copy($ws.Range("A26")) | Out-Null
$ws.Range("F3").Copy($ws.Range("F27")) | Out-Null

$ws.Range("A26").Value = "TrainScheduling_ltrailways_removeDepartment"
$ws.Range("C26").Value = "CompanyManagement.removeDepartment"
$ws.Range("A27").Value = "TrainScheduling_ltrailways_removeDepartment"
$ws.Range("C27").Value = "CompanyManagement.removeDepartment"

Write-Output "done"
